# Scheduled runner: refresh market price / profit figures on Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2698.9443
$ws.Range("I40").Value = 2198.5386
$ws.Range("K40").Value = 2198.5386
$ws.Range("M40").Value = -2023.5386
$ws.Range("H64").Value = 4250
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -5496
$ws.Range("H67").Value = 4250
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -6716
$ws.Range("H74").Value = 4699.95
$ws.Range("I74").Value = 6000.1113
$ws.Range("J74").Value = 3636.182
$ws.Range("K74").Value = 6000.1113
$ws.Range("L74").Value = 3636.182
$ws.Range("M74").Value = -5064.1113
$ws.Range("N74").Value = -5508.182
$ws.Range("H76").Value = 2885.7144
$ws.Range("I76").Value = 2700
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 2700
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -2385
$ws.Range("N76").Value = -4630
$ws.Range("H77").Value = 4699.95
$ws.Range("I77").Value = 6000.1113
$ws.Range("J77").Value = 3636.182
$ws.Range("K77").Value = 30000.5565
$ws.Range("L77").Value = 18180.91
$ws.Range("M77").Value = -25320.5565
$ws.Range("N77").Value = -27540.91
$ws.Range("H79").Value = 2885.7144
$ws.Range("I79").Value = 2700
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 2700
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -1608
$ws.Range("N79").Value = -6184
$ws.Range("H135").Value = 280.6842
$ws.Range("I135").Value = 197.70589
$ws.Range("J135").Value = 986
$ws.Range("K135").Value = 1779.35301
$ws.Range("L135").Value = 8874
$ws.Range("M135").Value = 755.64699
$ws.Range("N135").Value = -13944
$ws.Range("H141").Value = 4226.4165
$ws.Range("I141").Value = 3045.2307
$ws.Range("K141").Value = 9135.6921
$ws.Range("M141").Value = -3955.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8530.120999999999
$ws.Range("I32").Value = 5433.62
$ws.Range("K32").Value = 5433.62
$ws.Range("M32").Value = -5146.62
$ws.Range("H63").Value = 2410.0938
$ws.Range("I63").Value = 2391.434
$ws.Range("K63").Value = 2391.434
$ws.Range("M63").Value = -1705.434
$ws.Range("H66").Value = 2410.0938
$ws.Range("I66").Value = 2391.434
$ws.Range("K66").Value = 11957.17
$ws.Range("M66").Value = -8525.170000000002
$ws.Range("H132").Value = 1483.4828
$ws.Range("I132").Value = 1020.62
$ws.Range("J132").Value = 4376.375
$ws.Range("K132").Value = 3061.86
$ws.Range("L132").Value = 13129.125
$ws.Range("M132").Value = -531.8600000000001
$ws.Range("N132").Value = -18189.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2455.3635
$ws.Range("I105").Value = 2326.125
$ws.Range("J105").Value = 2800
$ws.Range("K105").Value = 2326.125
$ws.Range("L105").Value = 2800
$ws.Range("M105").Value = -579.125
$ws.Range("N105").Value = -6294
$ws.Range("H122").Value = 33779.5
$ws.Range("J122").Value = 33779.5
$ws.Range("L122").Value = 33779.5
$ws.Range("N122").Value = -43579.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2053.3333
$ws.Range("I62").Value = 2200
$ws.Range("J62").Value = 1980
$ws.Range("K62").Value = 2200
$ws.Range("L62").Value = 1980
$ws.Range("M62").Value = -1576
$ws.Range("N62").Value = -3228
$ws.Range("H65").Value = 2053.3333
$ws.Range("I65").Value = 2200
$ws.Range("J65").Value = 1980
$ws.Range("K65").Value = 11000
$ws.Range("L65").Value = 9900
$ws.Range("M65").Value = -7880
$ws.Range("N65").Value = -16140
$ws.Range("H132").Value = 1178.1282
$ws.Range("I132").Value = 766.73334
$ws.Range("J132").Value = 2549.4443
$ws.Range("K132").Value = 2300.20002
$ws.Range("L132").Value = 7648.3329
$ws.Range("M132").Value = 229.7999799999998
$ws.Range("N132").Value = -12708.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 453.5
$ws.Range("I97").Value = 570
$ws.Range("J97").Value = 447.3684
$ws.Range("K97").Value = 1710
$ws.Range("L97").Value = 1342.1052
$ws.Range("M97").Value = -1214
$ws.Range("N97").Value = -2334.1052
$ws.Range("H109").Value = 5556718.5
$ws.Range("I109").Value = 1745
$ws.Range("J109").Value = 16666666
$ws.Range("K109").Value = 5235
$ws.Range("L109").Value = 49999998
$ws.Range("M109").Value = -4195
$ws.Range("N109").Value = -50002078
$ws.Range("H110").Value = 6521
$ws.Range("I110").Value = 4900
$ws.Range("J110").Value = 6926.25
$ws.Range("K110").Value = 14700
$ws.Range("L110").Value = 20778.75
$ws.Range("M110").Value = -10610
$ws.Range("N110").Value = -28958.75
$ws.Range("H111").Value = 6078.3335
$ws.Range("I111").Value = 9000
$ws.Range("J111").Value = 5494
$ws.Range("K111").Value = 27000
$ws.Range("L111").Value = 16482
$ws.Range("M111").Value = -23933
$ws.Range("N111").Value = -22616
$ws.Range("H112").Value = 5838.75
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 5838.75
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 17516.25
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -19732.25
$ws.Range("H113").Value = 3052939.8
$ws.Range("I113").Value = 756262.0600000001
$ws.Range("J113").Value = 4608753.5
$ws.Range("K113").Value = 2268786.18
$ws.Range("L113").Value = 13826260.5
$ws.Range("M113").Value = -2266616.18
$ws.Range("N113").Value = -13830600.5
$ws.Range("H115").Value = 1601.421
$ws.Range("I115").Value = 1374
$ws.Range("K115").Value = 4122
$ws.Range("M115").Value = -2947
$ws.Range("H118").Value = 1853531.6
$ws.Range("I118").Value = 752.375
$ws.Range("J118").Value = 3970993.8
$ws.Range("K118").Value = 2257.125
$ws.Range("L118").Value = 11912981.4
$ws.Range("M118").Value = -1014.125
$ws.Range("N118").Value = -11915467.4
$ws.Range("H121").Value = 62506944
$ws.Range("I121").Value = 437.5
$ws.Range("J121").Value = 83342450
$ws.Range("K121").Value = 1312.5
$ws.Range("L121").Value = 250027350
$ws.Range("M121").Value = -2.5
$ws.Range("N121").Value = -250029970
$ws.Range("H122").Value = 2000715.2
$ws.Range("I122").Value = 661
$ws.Range("J122").Value = 5000796.5
$ws.Range("K122").Value = 5949
$ws.Range("L122").Value = 45007168.5
$ws.Range("M122").Value = -3499
$ws.Range("N122").Value = -45012068.5
$ws.Range("H131").Value = 88279.28999999999
$ws.Range("I131").Value = 450
$ws.Range("J131").Value = 117555.72
$ws.Range("K131").Value = 1350
$ws.Range("L131").Value = 352667.16
$ws.Range("M131").Value = 3690
$ws.Range("N131").Value = -362747.16

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4959.6665
$ws.Range("I70").Value = 5052
$ws.Range("J70").Value = 4775
$ws.Range("K70").Value = 5052
$ws.Range("L70").Value = 4775
$ws.Range("M70").Value = -4782
$ws.Range("N70").Value = -5315
$ws.Range("H73").Value = 4959.6665
$ws.Range("I73").Value = 5052
$ws.Range("J73").Value = 4775
$ws.Range("K73").Value = 5052
$ws.Range("L73").Value = 4775
$ws.Range("M73").Value = -4116
$ws.Range("N73").Value = -6647
